$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 619 (weekly data refresh:
# a new, most-recent observation is added at the top of this product's
# series, pushing the existing history down by two rows).
$ws.Rows.Item(619).EntireRow.Insert()
$ws.Rows.Item(619).EntireRow.Insert()

# New row 619: Primera
$ws.Range("A619").Value = 3
$ws.Range("B619").Value = "Femacal de La Calera"
$ws.Range("C619").Value = "Coquimbo"
$ws.Range("D619").Value = 44753
$ws.Range("E619").Value = 5
$ws.Range("F619").Value = 100112006
$ws.Range("G619").Value = "Repollo"
$ws.Range("H619").Value = "Crespo record"
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 2150
$ws.Range("K619").Value = 1100
$ws.Range("L619").Value = 1200
$ws.Range("M619").Value = 1144
$ws.Range("N619").Value = "`$/unidad"
$ws.Range("O619").Value = "Provincia de Quillota"
$ws.Range("P619").Value = 1144
$ws.Range("Q619").Value = 1
$ws.Range("R619").Value = "Hortaliza"

# New row 620: Segunda
$ws.Range("A620").Value = 3
$ws.Range("B620").Value = "Femacal de La Calera"
$ws.Range("C620").Value = "Coquimbo"
$ws.Range("D620").Value = 44753
$ws.Range("E620").Value = 5
$ws.Range("F620").Value = 100112006
$ws.Range("G620").Value = "Repollo"
$ws.Range("H620").Value = "Crespo record"
$ws.Range("I620").Value = "Segunda"
$ws.Range("J620").Value = 1100
$ws.Range("K620").Value = 900
$ws.Range("L620").Value = 900
$ws.Range("M620").Value = 900
$ws.Range("N620").Value = "`$/unidad"
$ws.Range("O620").Value = "Provincia de Quillota"
$ws.Range("P620").Value = 900
$ws.Range("Q620").Value = 1
$ws.Range("R620").Value = "Hortaliza"
